$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "iaest-measure:superficie-agricola-utilizada"
$ws.Range("B3").Value = "medida"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("B5").Clear()

$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"
